$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the newly-played season and its winner to the bottom of the table
$ws.Range("A32").Value = "22/23"
$ws.Range("B32").Value = "Manchester City"

# Scroll the window down a bit and move the selection the way it was left
# after the edit (cursor parked just below the newly added row)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A33").Select()
